$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update name of Giljastølen waste point (row 4) to 'Giljastølen Camping'
$ws.Range("D4").Value = "Giljastølen Camping"

# Update name of the FV286/Nevland waste point (row 13) to 'Nevland'
$ws.Range("D13").Value = "Nevland"

# Widen column D to fit the new text (35 characters)
$ws.Columns("D").ColumnWidth = 34.15

# Move the active selection to D23 (matches the saved view state)
$ws.Range("D23").Select()
